# Done with 383. Ransom note
# Adds a new row (row 14) for LeetCode problem 383 "Ransom Note" to the
# study-plan sheet, pushing nothing else around (rows 15/16 keep their
# existing placeholder content), and moves the active-cell selection to I15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 already has the exact formatting we need for the new row:
#   A -> blue fill ("No." style), D/E/I -> wrap text, H -> red status font.
# Copying just the A:I cell range (not the whole row) reproduces those
# per-cell styles on row 14 without touching columns beyond I.
$ws.Range("A12:I12").Copy($ws.Range("A14"))

$ws.Range("A14").Value = 383
$ws.Range("B14").Value = "Ransome Note"
$ws.Range("C14").Value = "String"
$ws.Range("D14").Value = "String, Hash table, Counting"
$ws.Range("E14").Value = "Dict"
$ws.Range("F14").Value = "Easy"
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = "✅"
$ws.Range("I14").Value = "Given O(n) sol and didn't see any other sol"

# Matches the new row height (45) used by the diff.
$ws.Rows.Item(14).RowHeight = 45

# The author's selection ended up on I15 after filling in the new row.
[void]$ws.Range("I15").Select()
